$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find the paragraph whose trimmed text equals $matchText (exact,
# ignoring the trailing paragraph mark) and insert a brand-new plain
# paragraph containing $newText immediately after it.
# ---------------------------------------------------------------------------
function Insert-ParagraphAfterMatch($doc, $matchText, $newText) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text.TrimEnd([char]13)
        if ($t -eq $matchText) {
            $p.Range.InsertParagraphAfter()
            $newP = $doc.Paragraphs.Item($i + 1)
            $newP.Range.Text = $newText
            return $true
        }
    }
    return $false
}

# ---------------------------------------------------------------------------
# 1. Professional summary: "21 years" -> "15+ years"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Senior Software Engineer with 21 years of experience",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Senior Software Engineer with 15+ years of experience", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. FLEEM bullet (Progressive Change Campaign Committee)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "• Conceived, architected, and engineered FLEEM web application using Twilio API for thousands of simultaneous phone calls",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "• Conceived, architected, and engineered FLEEM web application using Twilio API handling tens of thousands of calls using emulated predictive dialer for regulated political surveys", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Salsa Labs - geospatial/CRM bullet
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "• Maintained and extended comprehensive geospatial analysis and reporting tools for Java-based CRM system",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "• Maintained and extended comprehensive geospatial analysis and reporting tools for Java-based CRM system used by tens of thousands of users simultaneously", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4. Salsa Labs - mapping/visualization bullet
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "• Integrated mapping and visualization tools for political campaign data analysis",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "• Integrated mapping and visualization tools for political campaign data analysis interfacing with Government and Activism APIs", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5. Salsa Labs - new bullet after "Collaborated with political strategists..."
# ---------------------------------------------------------------------------
Insert-ParagraphAfterMatch $d `
    "• Collaborated with political strategists to translate geospatial requirements into technical solutions" `
    "• Handled billions of records with millions of columns in high-performance CRM system" | Out-Null

# ---------------------------------------------------------------------------
# 6. Praxis Project - new bullet after "Managed technology infrastructure..."
# ---------------------------------------------------------------------------
Insert-ParagraphAfterMatch $d `
    "• Managed technology infrastructure supporting community health initiatives across multiple countries" `
    "• Architected and developed 25 Drupal sites to integrate with membership databases, activism CRMs and government agencies, under guidelines from Kellogg Foundation and Robert Wood Johnson Foundation" | Out-Null

# ---------------------------------------------------------------------------
# 7. Lake Research Partners - new bullet after "Developed innovative approaches..."
# ---------------------------------------------------------------------------
Insert-ParagraphAfterMatch $d `
    "• Developed innovative approaches to visualizing demographic and market data for enhanced client understanding" `
    "• Trained staff on building Python tooling for report generation and analysis" | Out-Null

# ---------------------------------------------------------------------------
# 8. Replace the EDUCATION heading + its two sub-headings with a single new
#    plain bullet paragraph: "• Trained staff on PHP/MySQL for data analysis
#    and reporting systems"
# ---------------------------------------------------------------------------
$eduStart = -1
$eduEnd = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13)
    if ($t -eq "EDUCATION") {
        $eduStart = $i
    }
    if ($t -eq "Bachelor of Arts in Political Science - University of California, Berkeley") {
        $eduEnd = $i
    }
}

if ($eduStart -gt 0 -and $eduEnd -ge $eduStart) {
    $pStart = $d.Paragraphs.Item($eduStart)
    $pEnd = $d.Paragraphs.Item($eduEnd)
    $rngEdu = $d.Range($pStart.Range.Start, $pEnd.Range.End)
    $rngEdu.Delete()

    $pBefore = $d.Paragraphs.Item($eduStart - 1)
    $pBefore.Range.InsertParagraphAfter()
    $newEduP = $d.Paragraphs.Item($eduStart)
    $newEduP.Range.Text = "• Trained staff on PHP/MySQL for data analysis and reporting systems"
}

Write-Host "Done"
